$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellValue {
    param($row, $col, $val, $isText)
    $c = $ws.Cells.Item($row, $col)
    if ($isText) {
        $c.NumberFormat = "@"
        $c.Value2 = [string]$val
    } else {
        $c.Value2 = [double]$val
    }
}

function Set-CellStyle {
    param($row, $col, $donorRow, $donorCol)
    $ws.Cells.Item($donorRow, $donorCol).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial(-4122) | Out-Null
}

# --- Header text updates (Volume/Number and date range) ---
$ws.Cells.Item(8,1).Value2 = "Volume 31   Number  28"
$ws.Cells.Item(9,3).Value2 = "Report Covering the Week  7/8/2024  Through  7/14/2024"

# --- Crime data table updates ---
Set-CellValue 14 4 1 $false
Set-CellStyle 14 4 39 3
Set-CellValue 14 5 -100 $false
Set-CellStyle 14 5 39 11
Set-CellValue 14 10 7 $false
Set-CellValue 14 11 -71.428571428571 $false
Set-CellValue 14 14 -92.307692307692 $false
Set-CellValue 15 4 "0" $true
Set-CellStyle 15 4 39 1
Set-CellValue 15 5 "***.*" $true
Set-CellStyle 15 5 39 1
Set-CellValue 15 6 2 $false
Set-CellValue 15 8 0 $false
Set-CellValue 16 3 5 $false
Set-CellValue 16 4 7 $false
Set-CellValue 16 5 -28.571428571428 $false
Set-CellValue 16 6 21 $false
Set-CellValue 16 8 -8.695652173913 $false
Set-CellValue 16 9 140 $false
Set-CellValue 16 10 167 $false
Set-CellValue 16 11 -16.167664670658 $false
Set-CellValue 16 12 -4.761904761904 $false
Set-CellValue 16 13 3.703703703703 $false
Set-CellValue 16 14 -68.253968253968 $false
Set-CellValue 17 3 13 $false
Set-CellValue 17 4 12 $false
Set-CellValue 17 5 8.333333333333 $false
Set-CellValue 17 6 43 $false
Set-CellValue 17 7 58 $false
Set-CellValue 17 8 -25.862068965517 $false
Set-CellValue 17 9 253 $false
Set-CellValue 17 10 244 $false
Set-CellValue 17 11 3.688524590163 $false
Set-CellValue 17 12 12.444444444444 $false
Set-CellValue 17 13 22.815533980582 $false
Set-CellValue 17 14 -1.556420233463 $false
Set-CellValue 18 3 "0" $true
Set-CellStyle 18 3 39 1
Set-CellValue 18 5 -100 $false
Set-CellValue 18 6 28 $false
Set-CellValue 18 8 250 $false
Set-CellValue 18 10 98 $false
Set-CellValue 18 11 20.408163265306 $false
Set-CellValue 18 12 45.679012345679 $false
Set-CellValue 18 13 2.608695652173 $false
Set-CellValue 18 14 -77.307692307692 $false
Set-CellValue 19 3 8 $false
Set-CellValue 19 4 8 $false
Set-CellValue 19 5 0 $false
Set-CellValue 19 6 29 $false
Set-CellValue 19 7 39 $false
Set-CellValue 19 8 -25.641025641025 $false
Set-CellValue 19 9 222 $false
Set-CellValue 19 10 218 $false
Set-CellValue 19 11 1.834862385321 $false
Set-CellValue 19 12 -5.128205128205 $false
Set-CellValue 19 13 56.338028169014 $false
Set-CellValue 19 14 -3.056768558951 $false
Set-CellValue 20 3 "0" $true
Set-CellStyle 20 3 39 1
Set-CellValue 20 5 -100 $false
Set-CellValue 20 6 13 $false
Set-CellValue 20 7 25 $false
Set-CellValue 20 8 -48 $false
Set-CellValue 20 10 179 $false
Set-CellValue 20 11 -45.251396648044 $false
Set-CellValue 20 12 -21.6 $false
Set-CellValue 20 13 151.282051282051 $false
Set-CellValue 20 14 -54.838709677419 $false
Set-CellValue 21 3 26 $false
Set-CellValue 21 4 40 $false
Set-CellValue 21 5 -35 $false
Set-CellValue 21 6 136 $false
Set-CellValue 21 7 158 $false
Set-CellValue 21 8 -13.924050632911 $false
Set-CellValue 21 9 855 $false
Set-CellValue 21 10 929 $false
Set-CellValue 21 11 -7.965554359526 $false
Set-CellValue 21 12 2.888086642599 $false
Set-CellValue 21 13 31.741140215716 $false
Set-CellValue 21 14 -49.882766705744 $false
Set-CellValue 22 6 "0" $true
Set-CellStyle 22 6 39 1
Set-CellValue 22 7 2 $false
Set-CellValue 22 8 -100 $false
Set-CellValue 22 13 -30.769230769230 $false
Set-CellValue 23 4 "0" $true
Set-CellStyle 23 4 39 1
Set-CellValue 23 5 "***.*" $true
Set-CellStyle 23 5 39 1
Set-CellValue 24 3 16 $false
Set-CellValue 24 4 26 $false
Set-CellValue 24 5 -38.461538461538 $false
Set-CellValue 24 6 59 $false
Set-CellValue 24 7 75 $false
Set-CellValue 24 8 -21.333333333333 $false
Set-CellValue 24 9 417 $false
Set-CellValue 24 10 427 $false
Set-CellValue 24 11 -2.341920374707 $false
Set-CellValue 24 12 -7.126948775055 $false
Set-CellValue 24 13 39.464882943143 $false
Set-CellValue 25 3 2 $false
Set-CellValue 25 4 7 $false
Set-CellValue 25 5 -71.428571428571 $false
Set-CellValue 25 6 14 $false
Set-CellValue 25 7 21 $false
Set-CellValue 25 8 -33.333333333333 $false
Set-CellValue 25 9 111 $false
Set-CellValue 25 10 137 $false
Set-CellValue 25 11 -18.978102189781 $false
Set-CellValue 25 12 -43.367346938775 $false
Set-CellValue 26 3 12 $false
Set-CellValue 26 4 13 $false
Set-CellValue 26 5 -7.692307692307 $false
Set-CellValue 26 6 60 $false
Set-CellValue 26 7 45 $false
Set-CellValue 26 8 33.333333333333 $false
Set-CellValue 26 9 354 $false
Set-CellValue 26 10 274 $false
Set-CellValue 26 11 29.197080291970 $false
Set-CellValue 26 12 14.935064935064 $false
Set-CellValue 26 13 2.906976744186 $false
Set-CellValue 27 4 "0" $true
Set-CellStyle 27 4 39 1
Set-CellValue 27 5 "***.*" $true
Set-CellStyle 27 5 39 1
Set-CellValue 27 6 2 $false
Set-CellValue 27 7 4 $false
Set-CellValue 27 8 -50 $false
Set-CellValue 28 4 2 $false
Set-CellValue 28 5 -50 $false
Set-CellValue 28 6 12 $false
Set-CellValue 28 8 20 $false
Set-CellValue 28 9 71 $false
Set-CellValue 28 10 55 $false
Set-CellValue 28 11 29.090909090909 $false
Set-CellValue 28 12 14.516129032258 $false
Set-CellValue 29 3 1 $false
Set-CellStyle 29 3 39 3
Set-CellValue 29 4 "0" $true
Set-CellStyle 29 4 39 1
Set-CellValue 29 5 "***.*" $true
Set-CellStyle 29 5 39 1
Set-CellValue 29 6 2 $false
Set-CellValue 29 8 -50 $false
Set-CellValue 29 9 12 $false
Set-CellValue 29 11 -20 $false
Set-CellValue 29 12 -7.692307692307 $false
Set-CellValue 29 14 -77.358490566037 $false
Set-CellValue 30 3 1 $false
Set-CellStyle 30 3 39 3
Set-CellValue 30 4 "0" $true
Set-CellStyle 30 4 39 1
Set-CellValue 30 5 "***.*" $true
Set-CellStyle 30 5 39 1
Set-CellValue 30 6 2 $false
Set-CellValue 30 8 -33.333333333333 $false
Set-CellValue 30 9 11 $false
Set-CellValue 30 11 -15.384615384615 $false
Set-CellValue 30 12 0 $false
Set-CellValue 30 13 -8.333333333333 $false
Set-CellValue 30 14 -76.595744680851 $false

$excel.CutCopyMode = $false
